# Weekly data refresh: insert this week's two new price rows (row 161-162)
# for "1a (guarda)" / "2a (guarda)" Cebolla at Vega Monumental Concepción,
# pushing all the existing history rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 161 (shifts 161..245 -> 163..247).
$ws.Rows.Item(161).Resize(2).Insert()

# --- New row 161: "1a (guarda)" ---
$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 44460
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = 100112004
$ws.Range("G161").Value = "Cebolla"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "1a (guarda)"
$ws.Range("J161").Value = 600
$ws.Range("K161").Value = 6000
$ws.Range("L161").Value = 6500
$ws.Range("M161").Value = 6250
$ws.Range("N161").Value = "`$/malla 18 kilos"
$ws.Range("O161").Value = "Región de O'Higgins"
$ws.Range("P161").Value = 347
$ws.Range("Q161").Value = 18
$ws.Range("R161").Value = "Hortaliza"

# --- New row 162: "2a (guarda)" ---
$ws.Range("A162").Value = 11
$ws.Range("B162").Value = "Vega Monumental Concepción"
$ws.Range("C162").Value = "Bíobío"
$ws.Range("D162").Value = 44460
$ws.Range("E162").Value = 8
$ws.Range("F162").Value = 100112004
$ws.Range("G162").Value = "Cebolla"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "2a (guarda)"
$ws.Range("J162").Value = 300
$ws.Range("K162").Value = 5500
$ws.Range("L162").Value = 5500
$ws.Range("M162").Value = 5500
$ws.Range("N162").Value = "`$/malla 18 kilos"
$ws.Range("O162").Value = "Región de O'Higgins"
$ws.Range("P162").Value = 306
$ws.Range("Q162").Value = 18
$ws.Range("R162").Value = "Hortaliza"
